# Update column C (the "updated" date column) for rows 2-117:
# value 45188 (2023-09-19) -> 45189 (2023-09-20)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45188) {
        $cell.Value = 45189
    }
}
